$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the rows that changed in the repull
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F10").Value = -5
$ws.Range("F11").Value = -4
$ws.Range("F14").Value = -3
$ws.Range("F17").Value = -9
